$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1328.5714
$ws.Range("I32").Value = 1525
$ws.Range("J32").Value = 1066.6666
$ws.Range("K32").Value = 1525
$ws.Range("L32").Value = 1066.6666
$ws.Range("M32").Value = -1199
$ws.Range("N32").Value = -1718.6666

$ws.Range("H43").Value = 1622.1177
$ws.Range("I43").Value = 1382.3077
$ws.Range("J43").Value = 2401.5
$ws.Range("K43").Value = 1382.3077
$ws.Range("L43").Value = 2401.5
$ws.Range("M43").Value = -1313.3077
$ws.Range("N43").Value = -2539.5

$ws.Range("H74").Value = 2586.0356
$ws.Range("I74").Value = 2505.1904
$ws.Range("J74").Value = 2828.5715
$ws.Range("K74").Value = 2505.1904
$ws.Range("L74").Value = 2828.5715
$ws.Range("M74").Value = -1569.1904
$ws.Range("N74").Value = -4700.5715

$ws.Range("H77").Value = 2586.0356
$ws.Range("I77").Value = 2505.1904
$ws.Range("J77").Value = 2828.5715
$ws.Range("K77").Value = 12525.952
$ws.Range("L77").Value = 14142.8575
$ws.Range("M77").Value = -7845.951999999999
$ws.Range("N77").Value = -23502.8575

$ws.Range("H80").Value = 1171.2
$ws.Range("I80").Value = 781
$ws.Range("J80").Value = 1231.2307
$ws.Range("K80").Value = 2343
$ws.Range("L80").Value = 3693.6921
$ws.Range("M80").Value = -1345
$ws.Range("N80").Value = -5689.6921

$ws.Range("H83").Value = 1171.2
$ws.Range("I83").Value = 781
$ws.Range("J83").Value = 1231.2307
$ws.Range("K83").Value = 7029
$ws.Range("L83").Value = 11081.0763
$ws.Range("M83").Value = -2037
$ws.Range("N83").Value = -21065.0763

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 166.66667
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -84
$ws.Range("N4").Value = -332

$ws.Range("H6").Value = 24934
$ws.Range("J6").Value = 4800
$ws.Range("L6").Value = 4800
$ws.Range("N6").Value = -5146

$ws.Range("H9").Value = 21803.6
$ws.Range("I9").Value = 9000
$ws.Range("J9").Value = 25004.5
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 25004.5
$ws.Range("M9").Value = -8830
$ws.Range("N9").Value = -25344.5

$ws.Range("H20").Value = 21803.6
$ws.Range("I20").Value = 9000
$ws.Range("J20").Value = 25004.5
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 25004.5
$ws.Range("M20").Value = -8730
$ws.Range("N20").Value = -25544.5

$ws.Range("H23").Value = 72504.75
$ws.Range("J23").Value = 75003.5
$ws.Range("L23").Value = 75003.5
$ws.Range("N23").Value = -75521.5

$ws.Range("H32").Value = 16953558
$ws.Range("I32").Value = 17245826
$ws.Range("K32").Value = 17245826
$ws.Range("M32").Value = -17245539

$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -25976

$ws.Range("H55").Value = 24999.334
$ws.Range("J55").Value = 24999.334
$ws.Range("L55").Value = 24999.334
$ws.Range("N55").Value = -25629.334

$ws.Range("H80").Value = 42545
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 42545
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 42545
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -44541

$ws.Range("H83").Value = 42545
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 42545
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 127635
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -137619

$ws.Range("H132").Value = 6793.976
$ws.Range("I132").Value = 4589.472
$ws.Range("K132").Value = 13768.416
$ws.Range("M132").Value = -11238.416

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = $null

$ws.Range("H134").Value = 1513.2222
$ws.Range("I134").Value = 1419.28
$ws.Range("K134").Value = 4257.84
$ws.Range("M134").Value = -1722.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55.285713
$ws.Range("I7").Value = 33
$ws.Range("J7").Value = 68.454544
$ws.Range("K7").Value = 33
$ws.Range("L7").Value = 68.454544
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = -294.454544

$ws.Range("H16").Value = 3019.3333
$ws.Range("I16").Value = 1457.5
$ws.Range("J16").Value = 9266.666999999999
$ws.Range("K16").Value = 1457.5
$ws.Range("L16").Value = 9266.666999999999
$ws.Range("M16").Value = -1170.5
$ws.Range("N16").Value = -9840.666999999999

$ws.Range("H58").Value = 1049.8889
$ws.Range("I58").Value = 969.72095
$ws.Range("J58").Value = 1363.2727
$ws.Range("K58").Value = 969.72095
$ws.Range("L58").Value = 1363.2727
$ws.Range("M58").Value = -766.72095
$ws.Range("N58").Value = -1769.2727

$ws.Range("H62").Value = 3892.8838
$ws.Range("I62").Value = 4099.8647
$ws.Range("J62").Value = 2616.5
$ws.Range("K62").Value = 4099.8647
$ws.Range("L62").Value = 2616.5
$ws.Range("M62").Value = -3475.8647
$ws.Range("N62").Value = -3864.5

$ws.Range("H65").Value = 3892.8838
$ws.Range("I65").Value = 4099.8647
$ws.Range("J65").Value = 2616.5
$ws.Range("K65").Value = 20499.3235
$ws.Range("L65").Value = 13082.5
$ws.Range("M65").Value = -17379.3235
$ws.Range("N65").Value = -19322.5

$ws.Range("H105").Value = 5529.75
$ws.Range("I105").Value = 4039.6667
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 4039.6667
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -2292.6667
$ws.Range("N105").Value = -13494

$ws.Range("H113").Value = 3019.3333
$ws.Range("I113").Value = 1457.5
$ws.Range("J113").Value = 9266.666999999999
$ws.Range("K113").Value = 1457.5
$ws.Range("L113").Value = 9266.666999999999
$ws.Range("M113").Value = 712.5
$ws.Range("N113").Value = -13606.667

$ws.Range("H134").Value = 2525.6086
$ws.Range("I134").Value = 2347.8333
$ws.Range("K134").Value = 7043.499899999999
$ws.Range("M134").Value = -4508.499899999999

$ws.Range("H136").Value = 1049.8889
$ws.Range("I136").Value = 969.72095
$ws.Range("J136").Value = 1363.2727
$ws.Range("K136").Value = 2909.16285
$ws.Range("L136").Value = 4089.8181
$ws.Range("M136").Value = -359.1628500000002
$ws.Range("N136").Value = -9189.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4625
$ws.Range("I87").Value = 3500
$ws.Range("J87").Value = 8000
$ws.Range("K87").Value = 10500
$ws.Range("L87").Value = 24000
$ws.Range("M87").Value = -9252
$ws.Range("N87").Value = -26496

$ws.Range("H90").Value = 4625
$ws.Range("I90").Value = 3500
$ws.Range("J90").Value = 8000
$ws.Range("K90").Value = 31500
$ws.Range("L90").Value = 72000
$ws.Range("M90").Value = -25260
$ws.Range("N90").Value = -84480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 19500
$ws.Range("J39").Value = 19500
$ws.Range("L39").Value = 19500
$ws.Range("N39").Value = -20564

$ws.Range("H70").Value = 4638.737
$ws.Range("I70").Value = 4509.923
$ws.Range("J70").Value = 4917.8335
$ws.Range("K70").Value = 4509.923
$ws.Range("L70").Value = 4917.8335
$ws.Range("M70").Value = -4239.923
$ws.Range("N70").Value = -5457.8335

$ws.Range("H73").Value = 4638.737
$ws.Range("I73").Value = 4509.923
$ws.Range("J73").Value = 4917.8335
$ws.Range("K73").Value = 4509.923
$ws.Range("L73").Value = 4917.8335
$ws.Range("M73").Value = -3573.923
$ws.Range("N73").Value = -6789.8335

$ws.Range("H80").Value = 2536.2917
$ws.Range("I80").Value = 2235.5557
$ws.Range("J80").Value = 2716.7334
$ws.Range("K80").Value = 2235.5557
$ws.Range("L80").Value = 2716.7334
$ws.Range("M80").Value = -1237.5557
$ws.Range("N80").Value = -4712.7334

$ws.Range("H83").Value = 2536.2917
$ws.Range("I83").Value = 2235.5557
$ws.Range("J83").Value = 2716.7334
$ws.Range("K83").Value = 11177.7785
$ws.Range("L83").Value = 13583.667
$ws.Range("M83").Value = -6185.7785
$ws.Range("N83").Value = -23567.667

$ws.Range("H132").Value = 8084.952
$ws.Range("I132").Value = 12727.637
$ws.Range("K132").Value = 38182.911
$ws.Range("M132").Value = -35652.911

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2724.6365
$ws.Range("I46").Value = 3267.75
$ws.Range("J46").Value = 2414.2856
$ws.Range("K46").Value = 3267.75
$ws.Range("L46").Value = 2414.2856
$ws.Range("M46").Value = -3079.75
$ws.Range("N46").Value = -2790.2856

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5899.2085
$ws.Range("I136").Value = 13077
$ws.Range("K136").Value = 39231
$ws.Range("M136").Value = -36681
